# Generate Report for Handoff
# This script moves the two tracked localization files from "handed back"
# state to a fresh "ready for handoff" state: new source GUIDs, a reset
# status/date, and cleared handback/target bookkeeping columns.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "308f3b28-a034-47cb-9b43-3f727e5c105a"
$oldGuid2 = "6b083ed6-3a0c-4639-82f7-7375ea43d6fd"
$newGuid1 = "c6a4730e-947f-4e11-ba7d-89b71c45c3e2"
$newGuid2 = "ffff75181693-bf86-4d53-ae66-6896ec265cb4"

$newHash = "bf369d3b2f4a3bba97c0aa33ab7b3f274a7b8bcb"

$status = "Ready for handoff"
$emptyDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$ov2Url = $wsOverview.Hyperlinks.Item(1).Address
$ov3Url = $wsOverview.Hyperlinks.Item(2).Address

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("E2").Value = $status
$wsOverview.Range("F2").Value = $status
$wsOverview.Range("G2").Value = "2016-09-05 03:10:33"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-09-05 03:10:33"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B3").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ov2Url, "", "", "e2e\$newGuid1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ov3Url, "", "", "e2e\$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zh2Url = $wsZh.Hyperlinks.Item(1).Address
$zh3Url = $wsZh.Hyperlinks.Item(3).Address

$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("C2").Value = $status
$wsZh.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-05 03:10:28"
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $emptyDate

$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-05 03:10:28"
$wsZh.Range("I3").Hyperlinks.Delete()
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $emptyDate

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A3").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zh2Url, "", "", "$newGuid1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zh3Url, "", "", "$newGuid2.md")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$de2Url = $wsDe.Hyperlinks.Item(1).Address
$de3Url = $wsDe.Hyperlinks.Item(3).Address

$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("C2").Value = $status
$wsDe.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-05 03:10:33"
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $emptyDate

$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-05 03:10:33"
$wsDe.Range("I3").Hyperlinks.Delete()
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $emptyDate

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A3").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $de2Url, "", "", "$newGuid1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $de3Url, "", "", "$newGuid2.md")
